# "priority rule to be smarter"
# - Re-prioritise / expand the award (獎項) list on sheet 2: iPhone moves to
#   top spot (renamed "iPhone Xs"), and the big lump-sum prizes get split
#   into several graded tiers (2000元 x2, 1000元 x3, 600元 x5).
# - Selection moves from 工作表 (A1) to 獎項 (A2), which becomes the active tab.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # 工作表
$ws2 = $wb.Worksheets.Item(2)   # 獎項

# Reference cell that already carries the "numbered row" style used
# throughout the award table, so newly written rows pick up the same look.
$numberedStyle = $ws2.Range("A6").Style

# --- Rewrite the award column (B2:B6) with the new priority order ---
$ws2.Range("B2").Value = "iPhone Xs"
$ws2.Range("B3").Value = "3萬元"
$ws2.Range("B4").Value = "2萬元"
$ws2.Range("B5").Value = "PS4"
$ws2.Range("B6").Value = "5000元"

# --- Append the new, finer-grained prize tiers (rows 7-16) ---
$newAwards = @{
    7  = "2000元"
    8  = "2000元"
    9  = "1000元"
    10 = "1000元"
    11 = "1000元"
    12 = "600元"
    13 = "600元"
    14 = "600元"
    15 = "600元"
    16 = "600元"
}

for ($r = 7; $r -le 16; $r++) {
    $ws2.Range("A$r").Value = ($r - 1)
    $ws2.Range("A$r").Style = $numberedStyle
    $ws2.Range("B$r").Value = $newAwards[$r]
    $ws2.Range("B$r").Style = $numberedStyle
}

# --- Update selections / active tab ---
# 工作表's old selection (D5) moves on without staying the active sheet.
$ws1.Range("H13").Select() | Out-Null

# 獎項 becomes the active sheet, selection moves to D3.
$ws2.Activate() | Out-Null
$ws2.Range("D3").Select() | Out-Null
